$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50: phone column B50 changes from text "+919322612069" to numeric 919322612069
$ws.Cells.Item(50, 2).Value = 919322612069

# Row 51
$ws.Cells.Item(51, 1).Value = '2025-04-06 22:05:19'
$ws.Cells.Item(51, 2).Value = 919322612069
$ws.Cells.Item(51, 3).Value = 'text'
$ws.Cells.Item(51, 4).Value = 'hello'
$ws.Cells.Item(51, 5).Value = 'success'
$ws.Cells.Item(51, 6).Value = ''

# Row 52
$ws.Cells.Item(52, 1).Value = '2025-04-06 22:05:51'
$ws.Cells.Item(52, 2).Value = 91635348180
$ws.Cells.Item(52, 3).Value = 'text'
$ws.Cells.Item(52, 4).Value = 'hello'
$ws.Cells.Item(52, 5).Value = 'failed'
$ws.Cells.Item(52, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 53
$ws.Cells.Item(53, 1).Value = '2025-04-06 22:06:28'
$ws.Cells.Item(53, 2).Value = 919322612069
$ws.Cells.Item(53, 3).Value = 'text'
$ws.Cells.Item(53, 4).Value = 'hello'
$ws.Cells.Item(53, 5).Value = 'success'
$ws.Cells.Item(53, 6).Value = ''

# Row 54
$ws.Cells.Item(54, 1).Value = '2025-04-06 22:06:39'
$ws.Cells.Item(54, 2).Value = 916353481830
$ws.Cells.Item(54, 3).Value = 'text'
$ws.Cells.Item(54, 4).Value = 'hello'
$ws.Cells.Item(54, 5).Value = 'success'
$ws.Cells.Item(54, 6).Value = ''

# Row 55
$ws.Cells.Item(55, 1).Value = '2025-04-06 22:11:03'
$ws.Cells.Item(55, 2).Value = 919322612069
$ws.Cells.Item(55, 3).Value = 'text'
$ws.Cells.Item(55, 4).Value = 'yo'
$ws.Cells.Item(55, 5).Value = 'failed'
$ws.Cells.Item(55, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 56
$ws.Cells.Item(56, 1).Value = '2025-04-06 22:12:20'
$ws.Cells.Item(56, 2).Value = 919322612069
$ws.Cells.Item(56, 3).Value = 'text'
$ws.Cells.Item(56, 4).Value = 'yp'
$ws.Cells.Item(56, 5).Value = 'failed'
$ws.Cells.Item(56, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 57
$ws.Cells.Item(57, 1).Value = '2025-04-06 22:27:17'
$ws.Cells.Item(57, 2).Value = 918849958013
$ws.Cells.Item(57, 3).Value = 'text'
$ws.Cells.Item(57, 4).Value = 'cvfgh'
$ws.Cells.Item(57, 5).Value = 'success'
$ws.Cells.Item(57, 6).Value = ''

# Row 58
$ws.Cells.Item(58, 1).Value = '2025-04-06 22:37:05'
$ws.Cells.Item(58, 2).Value = 918849958013
$ws.Cells.Item(58, 3).Value = 'image'
$ws.Cells.Item(58, 4).Value = 'Image: harshan_attar.png'
$ws.Cells.Item(58, 5).Value = 'failed'
$ws.Cells.Item(58, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 59
$ws.Cells.Item(59, 1).Value = '2025-04-06 22:38:16'
$ws.Cells.Item(59, 2).Value = 918849958013
$ws.Cells.Item(59, 3).Value = 'image'
$ws.Cells.Item(59, 4).Value = 'Image: harshan_attar.png'
$ws.Cells.Item(59, 5).Value = 'failed'
$ws.Cells.Item(59, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 60
$ws.Cells.Item(60, 1).Value = '2025-04-06 22:50:09'
$ws.Cells.Item(60, 2).Value = 918849958013
$ws.Cells.Item(60, 3).Value = 'image'
$ws.Cells.Item(60, 4).Value = 'Image: WhatsApp_Image_2024-08-31_at_11.24.11_ae2e19bd-removebg-preview.png'
$ws.Cells.Item(60, 5).Value = 'failed'
$ws.Cells.Item(60, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 61
$ws.Cells.Item(61, 1).Value = '2025-04-06 22:54:54'
$ws.Cells.Item(61, 2).Value = 919322612069
$ws.Cells.Item(61, 3).Value = 'image'
$ws.Cells.Item(61, 4).Value = 'Image: 3853ec82-5d7e-4c3e-b266-2fca6f594a06.png'
$ws.Cells.Item(61, 5).Value = 'failed'
$ws.Cells.Item(61, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 62
$ws.Cells.Item(62, 1).Value = '2025-04-06 22:59:33'
$ws.Cells.Item(62, 2).Value = 918849958013
$ws.Cells.Item(62, 3).Value = 'image'
$ws.Cells.Item(62, 4).Value = 'Image: 65322a94-e7e5-4b4e-b7e7-f0eea51b01de.jpg'
$ws.Cells.Item(62, 5).Value = 'failed'
$ws.Cells.Item(62, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 63
$ws.Cells.Item(63, 1).Value = '2025-04-06 23:03:34'
$ws.Cells.Item(63, 2).Value = 918849958013
$ws.Cells.Item(63, 3).Value = 'image'
$ws.Cells.Item(63, 4).Value = 'Image: f41b2bcd-91ea-4c2b-8d38-af49c491fe4c.jpg'
$ws.Cells.Item(63, 5).Value = 'failed'
$ws.Cells.Item(63, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 64
$ws.Cells.Item(64, 1).Value = '2025-04-06 23:12:40'
$ws.Cells.Item(64, 2).Value = 918849958013
$ws.Cells.Item(64, 3).Value = 'image'
$ws.Cells.Item(64, 4).Value = 'Image: e30546a7-5841-4a13-a08d-7460ec2eb101.png'
$ws.Cells.Item(64, 5).Value = 'failed'
$ws.Cells.Item(64, 6).Value = 'Message: 
Stacktrace:
	GetHandleVerifier [0x00007FF68BA71F85+78133]
	GetHandleVerifier [0x00007FF68BA71FE0+78224]
	(No symbol) [0x00007FF68B8391BA]
	(No symbol) [0x00007FF68B88F19D]
	(No symbol) [0x00007FF68B88F44C]
	(No symbol) [0x00007FF68B8E23D7]
	(No symbol) [0x00007FF68B8B719F]
	(No symbol) [0x00007FF68B8DF21F]
	(No symbol) [0x00007FF68B8B6F33]
	(No symbol) [0x00007FF68B880358]
	(No symbol) [0x00007FF68B8810C3]
	GetHandleVerifier [0x00007FF68BD3BAAD+3001437]
	GetHandleVerifier [0x00007FF68BD35E92+2977858]
	GetHandleVerifier [0x00007FF68BD5499D+3103565]
	GetHandleVerifier [0x00007FF68BA8C81A+186826]
	GetHandleVerifier [0x00007FF68BA9442F+218591]
	GetHandleVerifier [0x00007FF68BA79DC4+110452]
	GetHandleVerifier [0x00007FF68BA79F72+110882]
	GetHandleVerifier [0x00007FF68BA603A9+5465]
	BaseThreadInitThunk [0x00007FFD0A3EE8D7+23]
	RtlUserThreadStart [0x00007FFD0B8DBF6C+44]
'

# Row 65
$ws.Cells.Item(65, 1).Value = '2025-04-06 23:17:51'
$ws.Cells.Item(65, 2).Value = 918849958013
$ws.Cells.Item(65, 3).Value = 'text'
$ws.Cells.Item(65, 4).Value = 'szxdcfbhnjmk'
$ws.Cells.Item(65, 5).Value = 'success'
$ws.Cells.Item(65, 6).Value = ''

# Row 66
$ws.Cells.Item(66, 1).Value = '2025-04-06 23:24:59'
$ws.Cells.Item(66, 2).Value = 919284374259
$ws.Cells.Item(66, 3).Value = 'text'
$ws.Cells.Item(66, 4).Value = 'sdfghjkl;'
$ws.Cells.Item(66, 5).Value = 'success'
$ws.Cells.Item(66, 6).Value = ''

# Row 67
$ws.Cells.Item(67, 1).Value = '2025-04-07 20:57:32'
$ws.Cells.Item(67, 2).Value = 919322612069
$ws.Cells.Item(67, 3).Value = 'text'
$ws.Cells.Item(67, 4).Value = 'hello'
$ws.Cells.Item(67, 5).Value = 'success'
$ws.Cells.Item(67, 6).Value = ''

# Row 68
$ws.Cells.Item(68, 1).Value = '2025-04-07 20:57:45'
$ws.Cells.Item(68, 2).Value = 916353481830
$ws.Cells.Item(68, 3).Value = 'text'
$ws.Cells.Item(68, 4).Value = 'hello'
$ws.Cells.Item(68, 5).Value = 'success'
$ws.Cells.Item(68, 6).Value = ''

# Row 69
$ws.Cells.Item(69, 1).Value = '2025-04-07 21:12:57'
$ws.Cells.Item(69, 2).Value = 918849958013
$ws.Cells.Item(69, 3).Value = 'text'
$ws.Cells.Item(69, 4).Value = 'Hello John Doe!'
$ws.Cells.Item(69, 5).Value = 'success'
$ws.Cells.Item(69, 6).Value = ''

# Row 70
$ws.Cells.Item(70, 1).Value = '2025-04-07 21:13:09'
$ws.Cells.Item(70, 2).Value = 919322612069
$ws.Cells.Item(70, 3).Value = 'text'
$ws.Cells.Item(70, 4).Value = 'Hello John Doe!'
$ws.Cells.Item(70, 5).Value = 'success'
$ws.Cells.Item(70, 6).Value = ''

# Row 71
$ws.Cells.Item(71, 1).Value = '2025-04-07 21:13:20'
$ws.Cells.Item(71, 2).Value = 917020811776
$ws.Cells.Item(71, 3).Value = 'text'
$ws.Cells.Item(71, 4).Value = 'Hello John Doe!'
$ws.Cells.Item(71, 5).Value = 'success'
$ws.Cells.Item(71, 6).Value = ''

# Row 72
$ws.Cells.Item(72, 1).Value = '2025-04-07 21:13:32'
$ws.Cells.Item(72, 2).Value = 919824237224
$ws.Cells.Item(72, 3).Value = 'text'
$ws.Cells.Item(72, 4).Value = 'Hello John Doe!'
$ws.Cells.Item(72, 5).Value = 'success'
$ws.Cells.Item(72, 6).Value = ''

# Row 73
$ws.Cells.Item(73, 1).Value = '2025-04-07 21:23:36'
$ws.Cells.Item(73, 2).Value = 918849958013
$ws.Cells.Item(73, 3).Value = 'text'
$ws.Cells.Item(73, 4).Value = 'Hello Yash Adagale!'
$ws.Cells.Item(73, 5).Value = 'success'
$ws.Cells.Item(73, 6).Value = ''

# Row 74
$ws.Cells.Item(74, 1).Value = '2025-04-07 21:23:47'
$ws.Cells.Item(74, 2).Value = 919322612069
$ws.Cells.Item(74, 3).Value = 'text'
$ws.Cells.Item(74, 4).Value = 'Hello Yash Adagale!'
$ws.Cells.Item(74, 5).Value = 'success'
$ws.Cells.Item(74, 6).Value = ''

# Row 75
$ws.Cells.Item(75, 1).Value = '2025-04-07 21:23:58'
$ws.Cells.Item(75, 2).Value = 919824237224
$ws.Cells.Item(75, 3).Value = 'text'
$ws.Cells.Item(75, 4).Value = 'Hello Yash Adagale!'
$ws.Cells.Item(75, 5).Value = 'success'
$ws.Cells.Item(75, 6).Value = ''

# Row 76
$ws.Cells.Item(76, 1).Value = '2025-04-07 21:27:23'
$ws.Cells.Item(76, 2).Value = 918849958013
$ws.Cells.Item(76, 3).Value = 'text'
$ws.Cells.Item(76, 4).Value = "'" + '1234567890'
$ws.Cells.Item(76, 5).Value = 'success'
$ws.Cells.Item(76, 6).Value = ''

# Row 77
$ws.Cells.Item(77, 1).Value = '2025-04-07 21:27:35'
$ws.Cells.Item(77, 2).Value = 919322612069
$ws.Cells.Item(77, 3).Value = 'text'
$ws.Cells.Item(77, 4).Value = "'" + '1234567890'
$ws.Cells.Item(77, 5).Value = 'success'
$ws.Cells.Item(77, 6).Value = ''

# Row 78
$ws.Cells.Item(78, 1).Value = '2025-04-07 21:27:47'
$ws.Cells.Item(78, 2).Value = "'" + '919824237224'
$ws.Cells.Item(78, 3).Value = 'text'
$ws.Cells.Item(78, 4).Value = "'" + '1234567890'
$ws.Cells.Item(78, 5).Value = 'success'
$ws.Cells.Item(78, 6).Value = ''

Write-Output "done"